$d = $word.ActiveDocument
$d.Content.Find.Execute("obvezna je vratiti se na rad dana 21. 8. 2024. godine.", $false, $false, $false, $false, $false, $true, 1, $false, "obvezna je vratiti se na rad dana {{ dpnr }}. {{ mjpnr }}. 2024. godine.", 2)
